# Adds five new country sheets (Norway, Poland, Portugal, Croatia, Greece)
# to the workbook, each cloned from the existing "Slovakia" template sheet,
# with the market name / ticket reference filled in on B2 / B4.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Give the template sheet a portrait page setup - new copies inherit it,
# matching the <pageSetup orientation="portrait"/> that shows up on every
# sheet from Slovakia onward in the final workbook.
$slovakia.PageSetup.Orientation = 1

# name, market-name cell (B2), ticket-ref cell (B4), whether B4 should be
# written before B2 (Poland is the one sheet where the ticket ref was
# entered first), and the final on-screen selection for that tab.
$countries = @(
    @{ Name = "Norway";   Market = "Norway Market";   Ticket = "NGC-2931/T3072/T3078/T3063"; TicketFirst = $false; Sel = "A1:XFD1048576" },
    @{ Name = "Poland";   Market = "Poland Market";   Ticket = "NGC-2920/T3038/T3106/T3121"; TicketFirst = $true;  Sel = "A1:XFD1048576" },
    @{ Name = "Portugal"; Market = "Portugal Market"; Ticket = "NGC-3479/T2410/T2430/T2461"; TicketFirst = $false; Sel = "E23" },
    @{ Name = "Croatia";  Market = "Croatia Market";  Ticket = "NGC-3139/T2418/T2474/T2488"; TicketFirst = $false; Sel = "A1:XFD1048576" },
    @{ Name = "Greece";   Market = "Greece Market";   Ticket = "NGC-4119/T3169/T3206/T3190"; TicketFirst = $false; Sel = "B4" }
)

$lastSheet = $null

foreach ($c in $countries) {
    $src = $wb.Worksheets.Item("Slovakia")
    $src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $c.Name

    if ($c.TicketFirst) {
        $newSheet.Range("B4").Value = $c.Ticket
        $newSheet.Range("B2").Value = $c.Market
    } else {
        $newSheet.Range("B2").Value = $c.Market
        $newSheet.Range("B4").Value = $c.Ticket
    }

    [void]$newSheet.Range($c.Sel).Select()
    $lastSheet = $newSheet
}

# The template sheet is no longer the active tab; its leftover selection
# becomes a full-sheet (column) selection like the other "parked" sheets.
[void]$slovakia.Range("A1:XFD1048576").Select()

# Greece (the last sheet created) stays the active tab/selection.
$lastSheet.Activate()
